$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.389.35'
$ws.Range("E2").Value = '  -3.85%  '
$ws.Range("D3").Value = '3.619.11'
$ws.Range("E3").Value = '  -4.25%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.39'
$ws.Range("E5").Value = '  -3.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '184.31'
$ws.Range("E6").Value = '  -2.23%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.612'
$ws.Range("E7").Value = '  -4.66%  '
$ws.Range("E8").Value = '  +0.20%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.675'
$ws.Range("E9").Value = '  -8.18%  '
$ws.Range("E10").Value = '  -12.52%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.73'
$ws.Range("E11").Value = '  -8.17%  '
$ws.Range("E12").Value = '  -15.93%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.98'
$ws.Range("E13").Value = '  -9.04%  '
$ws.Range("D14").Value = '4.190.71'
$ws.Range("E14").Value = '  -4.53%  '
$ws.Range("D15").Value = '3.614.91'
$ws.Range("E15").Value = '  -4.56%  '
$ws.Range("E16").Value = '  -0.70%  '
$ws.Range("D17").Value = '67.185.23'
$ws.Range("E17").Value = '  -3.87%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.43'
$ws.Range("E18").Value = '  -6.66%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.27'
$ws.Range("E19").Value = '  -6.56%  '
$ws.Range("E20").Value = '  -6.94%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '397.47'
$ws.Range("E21").Value = '  -5.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.32'
$ws.Range("E22").Value = '  -8.43%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '85.52'
$ws.Range("E23").Value = '  -5.72%  '
$ws.Range("E24").Value = '  -8.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.38'
$ws.Range("E25").Value = '  -5.89%  '
$ws.Range("E26").Value = '  -0.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.31'
$ws.Range("E27").Value = '  -9.51%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.64'
$ws.Range("E28").Value = '  -10.31%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.05'
$ws.Range("E29").Value = '  -6.70%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.26'
$ws.Range("E30").Value = '  -6.73%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.76'
$ws.Range("E31").Value = '  -10.83%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '66.03'
$ws.Range("E32").Value = '  +1.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.90'
$ws.Range("E33").Value = '  -6.32%  '
$ws.Range("B34").Value = 'Bittensor'
$ws.Range("C34").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '592.48'
$ws.Range("E34").Value = '  -4.24%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.113'
$ws.Range("E35").Value = '  -6.91%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '41.61'
$ws.Range("E36").Value = '  -7.91%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  -0.06%  '
$ws.Range("E39").Value = '  -8.92%  '
$ws.Range("E40").Value = '  -20.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.133'
$ws.Range("E41").Value = '  -5.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.78'
$ws.Range("E42").Value = '  -10.87%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0411'
$ws.Range("E43").Value = '  -8.94%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.42'
$ws.Range("E44").Value = '  -14.28%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '2.704.26'
$ws.Range("E45").Value = '  -3.98%  '
$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.130'
$ws.Range("E46").Value = '  -4.95%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.04'
$ws.Range("E47").Value = '  -6.64%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '138.97'
$ws.Range("E48").Value = '  -2.76%  '
$ws.Range("B49").Value = 'WEMIXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.55'
$ws.Range("E49").Value = '  -7.87%  '
$ws.Range("B50").Value = 'THORChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.37'
$ws.Range("E50").Value = '  -11.97%  '
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.58'
$ws.Range("E51").Value = '  -9.04%  '
